$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "Job ID" in F1
$ws.Range("F1").Value = "Job ID"

# Move the active selection to B6, matching the final selection state
$ws.Range("B6").Select()
